$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.769.45"
$ws.Range("E2").Value = "  +2.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.07"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.42"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4958"
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2962"
$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06808"
$ws.Range("E9").Value = "  +3.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.902.58"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.15"
$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07319"
$ws.Range("E12").Value = "  +2.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.19"
$ws.Range("E13").Value = "  +6.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.067"
$ws.Range("E14").Value = "  +4.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6747"
$ws.Range("E15").Value = "  +1.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.741.27"
$ws.Range("E16").Value = "  +2.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008003"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.31"
$ws.Range("E19").Value = "  +4.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.150.55"
$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.866"
$ws.Range("E22").Value = "  +2.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "177.72"
$ws.Range("E23").Value = "  +31.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.079"
$ws.Range("E24").Value = "  +8.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.312"
$ws.Range("E25").Value = "  +2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.54"
$ws.Range("E26").Value = "  +3.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.56"
$ws.Range("E27").Value = "  +10.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.935"
$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.376"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.357"
$ws.Range("E30").Value = "  +4.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08929"
$ws.Range("E31").Value = "  +3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  +2.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05261"
$ws.Range("E33").Value = "  +5.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7459"
$ws.Range("E34").Value = "  +5.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("E35").Value = "  +3.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.672"
$ws.Range("E36").Value = "  +0.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01893"
$ws.Range("E37").Value = "  +11.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.694"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.175"
$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9374"
$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4367"
$ws.Range("E41").Value = "  +4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.39"
$ws.Range("E42").Value = "  +3.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.841"
$ws.Range("E43").Value = "  -2.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.004"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.696"
$ws.Range("E45").Value = "  +3.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1372"
$ws.Range("E46").Value = "  +9.08%  "

$ws.Range("E47").Value = "  +3.08%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.595"
$ws.Range("E48").Value = "  +6.29%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.47"
$ws.Range("E49").Value = "  +3.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3894"
$ws.Range("E50").Value = "  +4.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.386"
$ws.Range("E51").Value = "  +3.58%  "
